$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: drop the trailing "no sequence diagram" remark (and its leading
# tab) from the OC-1 heading paragraph, leaving just the heading text.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "anmodOmRedigering`t//ingen sd da den er for simpel.",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "anmodOmRedigering",
    2) | Out-Null

# ---------------------------------------------------------------------------
# Change 2: the OC-2 system operation signature gains a second parameter.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "indtastNyKommentar(kommentar)",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "indtastNyKommentar(kommentar,koerselsid)",
    2) | Out-Null

# ---------------------------------------------------------------------------
# Change 3: the matching postcondition sentence spells out the same call.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "koersel.gemKommentar er blevet kaldt.",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "koersel.gemKommentar(kommentar,koerselsid) er blevet kaldt.",
    2) | Out-Null

# ---------------------------------------------------------------------------
# The hidden "_GoBack" bookmark (last-edit marker) that used to sit at the
# end of the OC-1 heading now belongs right before the closing ")" of the
# indtastNyKommentar(...) call that we just edited -- relocate it there.
# ---------------------------------------------------------------------------
$searchRange = $d.Content
$searchRange.Find.Execute(
    "indtastNyKommentar(kommentar,koerselsid)",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "",
    0) | Out-Null

$bookmarkPos = $searchRange.End - 1
$bookmarkRange = $d.Range($bookmarkPos, $bookmarkPos)
$d.Bookmarks.Add("_GoBack", $bookmarkRange) | Out-Null
